# Fruta / hortaliza, semanal
# Inserts 5 new weekly price rows (Limón, Macroferia Regional de Talca) right
# before the existing row 644, shifting the old rows 644-695 down to 649-700.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 644 (pushes old 644..695 to 649..700).
$ws.Rows.Item(644).Resize(5).Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$newRows = @(
        @(5, 'Macroferia Regional de Talca', 'Maule', 44461, 7, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '1a amarillo', 350, 4500, 4500, 4500, '$/malla 14 kilos', 'Provincia de Quillota', 321, 14),
        @(5, 'Macroferia Regional de Talca', 'Maule', 44461, 7, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '1a amarillo', 600, 4500, 4500, 4500, '$/malla 14 kilos', 'Región de O''Higgins', 321, 14),
        @(5, 'Macroferia Regional de Talca', 'Maule', 44461, 7, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '1a amarillo', 300, 5000, 5000, 5000, '$/malla 18 kilos', 'Provincia de Melipilla', 278, 18),
        @(5, 'Macroferia Regional de Talca', 'Maule', 44461, 7, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 250, 3500, 3500, 3500, '$/malla 14 kilos', 'Región de O''Higgins', 250, 14),
        @(5, 'Macroferia Regional de Talca', 'Maule', 44461, 7, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '3a amarillo', 200, 3000, 3000, 3000, '$/malla 14 kilos', 'Provincia de Quillota', 214, 14)
)

$startRow = 644
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowData = $newRows[$r]
    $rowNum = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value2 = $rowData[$c]
    }
}
